$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B-E keep their original text formatting (matches source data which
# is all stored as text/inlineStr, including numeric-looking price values).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.562.89"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.105.66"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.51"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.18"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.07"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.592.01"
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.59"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.108.41"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.993"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.96"
$ws.Range("E18").Value = "  +4.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.599.76"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.26"
$ws.Range("E20").Value = "  +7.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.97"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.13"
$ws.Range("E26").Value = "  -2.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.03"
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.23"
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.167"
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.106"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.38"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0480"
$ws.Range("E33").Value = "  +5.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.21"
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.10"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.292"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.29"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.62"
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.73"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.51"
$ws.Range("E47").Value = "  +4.59%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.067.15"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.937"
$ws.Range("E50").Value = "  +18.38%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0333"
$ws.Range("E51").Value = "  +3.94%  "
